# Atualização de bases das ligas, do dia: 17-06-2024 às 21:10
# Swap the match-data (columns B:AD) between pairs of adjacent rows,
# leaving column A (the running row id) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AD$rowA")
    $rangeB = $ws.Range("B$rowB`:AD$rowB")

    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

Swap-RowData 170 171
Swap-RowData 200 201
Swap-RowData 213 214
